$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.680.26"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.695.60"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'2.45"
$ws.Range("E4").Value = "  +30.93%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'228.82"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("D7").Value = "'651.93"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'0.438"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("D9").Value = "'1.15"
$ws.Range("E9").Value = "  +9.62%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "3.693.77"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'47.60"
$ws.Range("E12").Value = "  +7.59%  "
$ws.Range("D13").Value = "'0.210"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").Value = "'0.0000300"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'6.63"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "4.386.35"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "96.394.90"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'8.89"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "3.685.75"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'19.59"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").Value = "'12.89"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +9.28%  "
$ws.Range("D23").Value = "'532.68"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").Value = "'3.32"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  +43.53%  "
$ws.Range("D26").Value = "'121.28"
$ws.Range("E26").Value = "  +20.53%  "
$ws.Range("D27").Value = "'0.0000209"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "'6.82"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "3.895.16"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'12.96"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "'13.31"
$ws.Range("E31").Value = "  +9.97%  "
$ws.Range("D32").Value = "'2.99"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'0.186"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "'33.28"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.610"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "'606.68"
$ws.Range("E39").Value = "  -7.19%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'8.42"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").Value = "'7.14"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").Value = "'0.512"
$ws.Range("E43").Value = "  +20.06%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.163"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0503"
$ws.Range("E45").Value = "  +12.16%  "
$ws.Range("D46").Value = "'40.50"
$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "'0.970"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "'9.00"
$ws.Range("E49").Value = "  +6.65%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").Value = "'23.54"
$ws.Range("E51").Value = "  -0.10%  "
